# update AT0029 download template
# Refresh the rolling-date values on the "Example" settings sheet and
# leave the selection where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# targetMonth (E2): 2024-10 -> 2025-04
$ws.Range("E2").Value = 202504

# startDate (K2): 2024-10-04 -> 2025-03-04
$ws.Range("K2").Value = 20250304

# endDate (L2): 2024-11-03 -> 2025-04-03
$ws.Range("L2").Value = 20250403

# Move the active selection to J20, matching the saved view state.
$ws.Range("J20").Select()
